# Adding CSV, Excel and property reader.
# Rebuild sheet1 as a 6-column test-data table (aut/username/password/browser/
# headless/Pass-or-fail) with per-row hyperlinks on the credential columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- wipe old content / formatting / hyperlinks -----------------------------
$ws.Cells.ClearFormats()
$ws.Hyperlinks.Delete()
$ws.Cells.ClearContents()

# --- header row ---------------------------------------------------------
$ws.Range("A1").Value = "aut"
$ws.Range("B1").Value = "username"
$ws.Range("C1").Value = "password"
$ws.Range("D1").Value = "browser"
$ws.Range("E1").Value = "headless"
$ws.Range("F1").Value = "Pass or fail"

# --- data rows ------------------------------------------------------------
$ws.Range("A2").Value = "https://ui.cogmento.com"
$ws.Range("B2").Value = "Veercraig@gmail.com"
$ws.Range("C2").Value = "pD@zZVvh7pVWJbS"
$ws.Range("D2").Value = "chrome"
$ws.Range("E2").Value = $true
$ws.Range("F2").Value = $true

$ws.Range("A3").Value = "https://ui.cogmento.com"
$ws.Range("B3").Value = "Veercraig1@gmail.com"
$ws.Range("C3").Value = "pD@zZVvh7pVWJbS"
$ws.Range("D3").Value = "chrome"
$ws.Range("E3").Value = $true
$ws.Range("F3").Value = $false

$ws.Range("A4").Value = "https://ui.cogmento.com"
$ws.Range("B4").Value = "Veercraig@gmail.com"
$ws.Range("C4").Value = "pD@zZVvh7pVWJbS12"
$ws.Range("D4").Value = "chrome"
$ws.Range("E4").Value = $true
$ws.Range("F4").Value = $false

$ws.Range("A5").Value = "https://ui.cogmento.com"
$ws.Range("C5").Value = "pD@zZVvh7pVWJbS23"
$ws.Range("B5").Value = "Veercraig11@gmail.com"
$ws.Range("D5").Value = "chrome"
$ws.Range("E5").Value = $true
$ws.Range("F5").Value = $false

$ws.Range("A6").Value = "https://ui.cogmento.com"
$ws.Range("B6").Value = "Veercraig@gmail.com"
$ws.Range("C6").Value = "pD@zZVvh7pVWJbS"
$ws.Range("D6").Value = "chrome"
$ws.Range("E6").Value = $true
$ws.Range("F6").Value = $true

# --- column widths (approximate best-fit) ----------------------------------
$ws.Columns.Item(1).ColumnWidth = 21.830729166666668
$ws.Columns.Item(2).ColumnWidth = 18.276041666666668
$ws.Columns.Item(3).ColumnWidth = 16.498697916666668
$ws.Columns.Item(4).ColumnWidth = 6.830729166666667
$ws.Columns.Item(5).ColumnWidth = 7.166666666666667
$ws.Columns.Item(6).ColumnWidth = 8.944010416666666

# --- hyperlinks (auto-linked URL / mailto cells), in authoring order -------
$ws.Hyperlinks.Add($ws.Range("A2"), "https://ui.cogmento.com")
$ws.Hyperlinks.Add($ws.Range("A3"), "https://ui.cogmento.com")
$ws.Hyperlinks.Add($ws.Range("A4"), "https://ui.cogmento.com")
$ws.Hyperlinks.Add($ws.Range("A5"), "https://ui.cogmento.com")
$ws.Hyperlinks.Add($ws.Range("A6"), "https://ui.cogmento.com")
$ws.Hyperlinks.Add($ws.Range("B3"), "mailto:Veercraig1@gmail.com")
$ws.Hyperlinks.Add($ws.Range("C4"), "mailto:pD@zZVvh7pVWJbS12")
$ws.Hyperlinks.Add($ws.Range("C5"), "mailto:pD@zZVvh7pVWJbS23")
$ws.Hyperlinks.Add($ws.Range("B5"), "mailto:Veercraig11@gmail.com")

# --- selection back to A1 (matches a freshly laid-out sheet) ---------------
$null = $ws.Range("A1").Select()
